$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 963
$ws.Range("C3").Value = 1913
$ws.Range("D3").Value = 3855
$ws.Range("E3").Value = 10300
$ws.Range("F3").Value = 11600
$ws.Range("G3").Value = 12300

$ws.Range("B8").Value = 8246
$ws.Range("C8").Value = 14500
$ws.Range("D8").Value = 27600
$ws.Range("E8").Value = 55800
$ws.Range("F8").Value = 98500
$ws.Range("G8").Value = 126000

$ws.Range("B13").Value = 7211
$ws.Range("C13").Value = 12300
$ws.Range("D13").Value = 15800
$ws.Range("E13").Value = 18300
$ws.Range("F13").Value = 19100
$ws.Range("G13").Value = 18800

$ws.Range("B18").Value = 119000
$ws.Range("C18").Value = 188000
$ws.Range("D18").Value = 278000
$ws.Range("E18").Value = 357000
$ws.Range("F18").Value = 401000
$ws.Range("G18").Value = 468000

$ws.Range("B23").Value = 3190
$ws.Range("C23").Value = 4179
$ws.Range("D23").Value = 7627
$ws.Range("E23").Value = 7827
$ws.Range("F23").Value = 8255
$ws.Range("G23").Value = 13400

$ws.Range("B28").Value = 125000
$ws.Range("C28").Value = 135000
$ws.Range("D28").Value = 226000
$ws.Range("E28").Value = 256000
$ws.Range("F28").Value = 279000
$ws.Range("G28").Value = 253000

$ws.Range("B33").Value = 7501
$ws.Range("C33").Value = 10100
$ws.Range("D33").Value = 11700
$ws.Range("E33").Value = 12900
$ws.Range("F33").Value = 13600
$ws.Range("G33").Value = 14100

$ws.Range("B38").Value = 115000
$ws.Range("C38").Value = 189000
$ws.Range("D38").Value = 257000
$ws.Range("E38").Value = 305000
$ws.Range("F38").Value = 356000
$ws.Range("G38").Value = 418000
